# Update crypto price/volume data per Sun Feb 18 12:55:59 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.666.44'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '2.797.64'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = "'353.46"
$ws.Range("E5").Value = '  -1.34%  '
$ws.Range("D6").Value = "'111.11"
$ws.Range("E6").Value = '  +1.18%  '
$ws.Range("D7").Value = "'0.554"
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = "'0.629"
$ws.Range("E9").Value = '  +6.90%  '
$ws.Range("D10").Value = "'40.14"
$ws.Range("E10").Value = '  +0.86%  '
$ws.Range("E11").Value = '  -2.68%  '
$ws.Range("D12").Value = "'0.0838"
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("D13").Value = "'19.96"
$ws.Range("E13").Value = '  +0.65%  '
$ws.Range("D14").Value = "'7.74"
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").Value = '3.242.54'
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").Value = '2.811.00'
$ws.Range("E16").Value = '  +1.90%  '
$ws.Range("D17").Value = "'0.946"
$ws.Range("E17").Value = '  +1.79%  '
$ws.Range("D18").Value = '51.642.78'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = "'7.59"
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("B20").Value = 'ImmutableX'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D20").Value = "'3.19"
$ws.Range("E20").Value = '  +3.00%  '
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").Value = "'13.78"
$ws.Range("E21").Value = '  +4.07%  '
$ws.Range("D22").Value = '0.0₃0970'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = "'70.24"
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = "'266.96"
$ws.Range("E24").Value = '  -0.56%  '
$ws.Range("D25").Value = "'2.76"
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").Value = "'26.09"
$ws.Range("E27").Value = '  -1.10%  '
$ws.Range("E28").Value = '  -3.21%  '
$ws.Range("D29").Value = "'38.85"
$ws.Range("E29").Value = '  +9.26%  '
$ws.Range("D30").Value = "'10.33"
$ws.Range("E30").Value = '  +1.31%  '
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("D32").Value = "'52.56"
$ws.Range("E32").Value = '  +0.88%  '
$ws.Range("D33").Value = "'6.10"
$ws.Range("E33").Value = '  -1.66%  '
$ws.Range("D34").Value = "'0.0886"
$ws.Range("E34").Value = '  +5.53%  '
$ws.Range("D35").Value = "'5.59"
$ws.Range("E35").Value = '  +7.27%  '
$ws.Range("D36").Value = "'0.0444"
$ws.Range("E36").Value = '  +0.30%  '
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("D38").Value = "'18.82"
$ws.Range("E38").Value = '  -0.32%  '
$ws.Range("D39").Value = "'3.16"
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("E40").Value = '  +2.30%  '
$ws.Range("E41").Value = '  +0.56%  '
$ws.Range("D42").Value = "'2.49"
$ws.Range("E42").Value = '  -0.91%  '
$ws.Range("E43").Value = '  +1.73%  '
$ws.Range("D44").Value = "'121.17"
$ws.Range("E44").Value = '  +1.35%  '
$ws.Range("D45").Value = "'21.89"
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("D46").Value = "'2.46"
$ws.Range("E46").Value = '  +6.03%  '
$ws.Range("D47").Value = "'3.40"
$ws.Range("E47").Value = '  +4.20%  '
$ws.Range("D48").Value = '2.103.07'
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("D49").Value = "'0.952"
$ws.Range("E49").Value = '  +1.61%  '
$ws.Range("D50").Value = "'5.45"
$ws.Range("E50").Value = '  -1.87%  '
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").Value = "'1.35"
$ws.Range("E51").Value = '  +5.95%  '
